$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at G, shifting the old G (d=7) -> H and old H (d=10) -> I.
$ws.Columns("G").Insert()

# New header for the inserted column.
$ws.Range("G1").Value = "d=6"

# New values for the inserted column (rows 2-6).
$ws.Range("G2").Value = 97.96335434236535
$ws.Range("G3").Value = 98.10979447162063
$ws.Range("G4").Value = 97.95403931237809
$ws.Range("G5").Value = 97.97398492052123
$ws.Range("G6").Value = 98.01147503546618
